$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

# Rows 101 and 102 were previously blank placeholder rows (with only the
# shared "" formula in column E). Fill them in with two new logged work
# sessions (both categorized as "Coding"), continuing the time log.
$ws1.Range("A101").Value2 = 41939
$ws1.Range("B101").Value2 = 0.92222222222222217
$ws1.Range("C101").Value2 = 0.99236111111111114
$ws1.Range("D101").Value2 = 0
$ws1.Range("F101").Value2 = "Coding"

$ws1.Range("A102").Value2 = 41940
$ws1.Range("B102").Value2 = 0.52847222222222223
$ws1.Range("C102").Value2 = 0.57638888888888895
$ws1.Range("D102").Value2 = 0
$ws1.Range("F102").Value2 = "Coding"

# Push the remaining blank placeholder rows (old rows 103 and 104, the
# trailing blank formula row and the Total Time row) further down, giving
# the log room to grow - 18 new blank rows are inserted above them.
$ws1.Range("A103:A120").EntireRow.Insert() | Out-Null

# Update the view to match where the user ended up editing.
$ws1.Range("B121").Select() | Out-Null
